$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 and J1, matching the style of the existing header cells
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in column I (always 1) and column J (mirrors column H) for data rows 2-25
for ($row = 2; $row -le 25; $row++) {
    $hValue = $ws.Cells.Item($row, 8).Value2
    $ws.Cells.Item($row, 9).Value = 1
    $ws.Cells.Item($row, 10).Value = $hValue
}
